# Add the Spanish sitewide-search subsite sheet and tidy up the data providers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Sheets.Item(1)

# --- Create the new worksheet after the existing one ------------------------
$ws2 = $wb.Worksheets.Add($null, $wb.Sheets.Item($wb.Sheets.Count))
$ws2.Name = "SitewideSearchEs"

# Copy the bold header formatting from the English sheet's header row.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# --- Header row ---------------------------------------------------------
$ws2.Range("A1").Value = "CancerTerm"
$ws2.Range("B1").Value = "BestBet"
$ws2.Range("C1").Value = "Definition"

# --- Data (written column-by-column to mirror the authoring order) ------
$ws2.Range("B2").Value = "dolor"
$ws2.Range("B3").Value = "tabaco"
$ws2.Range("B4").Value = "linfoma"

$ws2.Range("C2").Value = "cáncer"
$ws2.Range("C3").Value = "glioma"
$ws2.Range("C4").Value = "cáncer de hígado"

$ws2.Range("A8").Value = "A33"

$ws2.Range("A2").Value = "safingol"
$ws2.Range("A3").Value = "dalteparina sódica"
$ws2.Range("A4").Value = "macrófago"
$ws2.Range("A5").Value = "dermis"
$ws2.Range("A6").Value = "TAC-101"
$ws2.Range("A7").Value = "lanolina"

# --- Column widths (best-fit, matching the English sheet's columns) -----
$ws2.Columns("A").ColumnWidth = 16.666666666666668
$ws2.Columns("B").ColumnWidth = 25.166666666666668
$ws2.Columns("C").ColumnWidth = 26.333333333333336

# --- Selection / active-sheet bookkeeping --------------------------------
[void]$ws1.Range("A9").Select()
[void]$ws2.Range("A9").Select()
[void]$ws2.Activate()
